# Update "4D Box" results sheet with the 2/7/2025 (Wed) draw.
# The newest draw already occupies row 2; this edit inserts a duplicate of
# the previous top entry (29/6/2025) as a new row 4, pushing the six rows
# below it down by one (through row 10), and relabels row 3 to read
# 2/7/2025 (Wed) so it lines up with the newest result in row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 2/7/2025 (Wed)
$ws.Range("A3").Value = "2/7/2025 (Wed)"
$ws.Range("B3").Value = "6 5 2 6`n4 0 7 9`n0 4 6 5`n1 3 0 8"
$ws.Range("C3").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

# Row 4: 29/6/2025 (Sun)
$ws.Range("A4").Value = "29/6/2025 (Sun)"
$ws.Range("B4").Value = "6 5 2 6`n4 0 7 9`n0 4 6 5`n1 3 0 8"
$ws.Range("C4").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

# Row 5: 28/6/2025 (Sat)
$ws.Range("A5").Value = "28/6/2025 (Sat)"
$ws.Range("B5").Value = "3 4 6 0`n4 9 3 6`n1 5 2 7`n0 0 4 8"
$ws.Range("C5").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

# Row 6: 25/6/2025 (Wed)
$ws.Range("A6").Value = "25/6/2025 (Wed)"
$ws.Range("B6").Value = "2 0 5 7`n6 2 8 8`n9 1 3 0`n7 6 1 4"
$ws.Range("C6").Value = "✅ Direct: 12/4302 (0.28%)`n✅ iBet: 12/226 (5.31%)"

# Row 7: 22/6/2025 (Sun)
$ws.Range("A7").Value = "22/6/2025 (Sun)"
$ws.Range("B7").Value = "4 1 3 7`n6 2 5 4`n0 4 2 8`n9 5 6 3"
$ws.Range("C7").Value = "✅ Direct: 11/4144 (0.27%)`n✅ iBet: 11/222 (4.95%)"

# Row 8: 21/6/2025 (Sat)
$ws.Range("A8").Value = "21/6/2025 (Sat)"
$ws.Range("B8").Value = "4 6 1 8`n8 1 0 4`n1 5 7 7`n0 2 9 3"
$ws.Range("C8").Value = "✅ Direct: 13/3814 (0.34%)`n✅ iBet: 13/208 (6.25%)"

# Row 9: 18/6/2025 (Wed)
$ws.Range("A9").Value = "18/6/2025 (Wed)"
$ws.Range("B9").Value = "2 1 2 1`n3 2 4 7`n0 5 6 9`n5 3 3 8"
$ws.Range("C9").Value = "✅ Direct: 9/3416 (0.26%)`n✅ iBet: 9/188 (4.79%)"

# Row 10: 15/6/2025 (Sun)
$ws.Range("A10").Value = "15/6/2025 (Sun)"
$ws.Range("B10").Value = "2 9 3 2`n0 2 8 9`n8 5 2 5`n6 7 4 1"
$ws.Range("C10").Value = "✅ Direct: 12/3547 (0.34%)`n✅ iBet: 12/195 (6.15%)"

# Row 13 gains a blank, formatted C13 cell (matching B13's style) to keep
# the placeholder block rectangular.
$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# A new blank, formatted row 39 is appended below row 38, extending the
# placeholder block by one row (matches the dimension growing to C39).
$ws.Range("B38").Copy()
$ws.Range("B39").PasteSpecial(-4122)

$excel.CutCopyMode = $false

